$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the extracted request data in column B, keeping column A labels intact
$ws.Range("B1").Value = "Moringa"
$ws.Range("B2").Value = "moringa@email.com"
$ws.Range("B3").Value = "Registration"
$ws.Range("B4").Value = "Patrick"
$ws.Range("B5").Value = 713636981

# Update the current selection to match the final state of the sheet
$ws.Range("D9").Select()

$wb.Save()
